$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Imputaciones semana")
$ws.Activate()

$ws.Range("A2:G2").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)

$ws.Range("A5").Value = 16
$ws.Range("B5").Value = 20.75
$ws.Range("C5").Value = 15
$ws.Range("D5").Value = 11.75
$ws.Range("E5").Value = 21
$ws.Range("F5").Value = 13.5
$ws.Range("G5").Value = 15.1
